# Natmi following Dr Hou advice
# Rebuild the LR-pair data table: recompute values for the existing
# "ECs -> Fndc5/Itgb5 -> {ECs,FAPs,sCs}" rows and add the new
# "FAPs -> Fndc5/Itgb5 -> {ECs,FAPs,sCs}" rows that Dr Hou asked for.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T for rows 2..7 (Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric score columns).
$data = @(
    @("ECs",  "Fndc5", "Itgb5", "ECs",  2, 0.6666666666666666, 0.116649, 0.349947, 0.186245347817965, 0.186245347817965, 3, 1, 6.126464333333334,  18.379393,  0.1081098818071741, 0.1081098818071741, 0.7146459380190001, 6.431813442171,     0.02013496253973623, 0.02013496253973623),
    @("ECs",  "Fndc5", "Itgb5", "FAPs", 2, 0.6666666666666666, 0.116649, 0.349947, 0.186245347817965, 0.186245347817965, 3, 1, 35.73736866666667,  107.212106, 0.6306349784216607, 0.6306349784216608, 4.168728317598,     37.518554858382,    0.117452830902317,   0.1174528309023171),
    @("ECs",  "Fndc5", "Itgb5", "sCs",  2, 0.6666666666666666, 0.116649, 0.349947, 0.186245347817965, 0.186245347817965, 3, 1, 14.80503233333333,  44.415097,  0.2612551397711651, 0.2612551397711651, 1.726992216651,     15.542929949859,    0.0486575543759117,  0.0486575543759117),
    @("FAPs", "Fndc5", "Itgb5", "ECs",  3, 1,                  0.50967,  1.52901,  0.813754652182035, 0.813754652182035, 3, 1, 6.126464333333334,  18.379393,  0.1081098818071741, 0.1081098818071741, 3.12247507677,      28.10227569093,     0.08797491926743789, 0.08797491926743789),
    @("FAPs", "Fndc5", "Itgb5", "FAPs", 3, 1,                  0.50967,  1.52901,  0.813754652182035, 0.813754652182035, 3, 1, 35.73736866666667,  107.212106, 0.6306349784216607, 0.6306349784216608, 18.21426468834,     163.92838219506,    0.5131821475193437,  0.5131821475193438),
    @("FAPs", "Fndc5", "Itgb5", "sCs",  3, 1,                  0.50967,  1.52901,  0.813754652182035, 0.813754652182035, 3, 1, 14.80503233333333,  44.415097,  0.2612551397711651, 0.2612551397711651, 7.54568082933,      67.91112746397,     0.2125975853952534,  0.2125975853952534)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
